$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the data range A2:C13 alphabetically by Block Name (column A), ascending
$rng = $ws.Range("A2:C13")
$sortField = $ws.Range("A2:A13")
$rng.Sort($sortField, 1)

# Clear the explicit style (s="1") from all cells in A1:C13 so they revert to default style
$ws.Range("A1:C13").Style = "Normal"
